# Auto-generated from diff: updates market-price / leve-profit figures
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 418.55
$ws.Range("I15").Value = 418.55
$ws.Range("K15").Value = 1255.65
$ws.Range("M15").Value = -1086.65
$ws.Range("H32").Value = 1053
$ws.Range("I32").Value = 1348.3334
$ws.Range("J32").Value = 610
$ws.Range("K32").Value = 1348.3334
$ws.Range("L32").Value = 610
$ws.Range("M32").Value = -1022.3334
$ws.Range("N32").Value = -1262
$ws.Range("H98").Value = 19601.182
$ws.Range("I98").Value = 23049.762
$ws.Range("J98").Value = 1975.1111
$ws.Range("K98").Value = 23049.762
$ws.Range("L98").Value = 1975.1111
$ws.Range("M98").Value = -21551.762
$ws.Range("N98").Value = -4971.1111
$ws.Range("H122").Value = 19601.182
$ws.Range("I122").Value = 23049.762
$ws.Range("J122").Value = 1975.1111
$ws.Range("K122").Value = 69149.28599999999
$ws.Range("L122").Value = 5925.3333
$ws.Range("M122").Value = -66699.28599999999
$ws.Range("N122").Value = -10825.3333
$ws.Range("H135").Value = 978.85297
$ws.Range("I135").Value = 131.2963
$ws.Range("J135").Value = 4248
$ws.Range("K135").Value = 1181.6667
$ws.Range("L135").Value = 38232
$ws.Range("M135").Value = 1353.3333
$ws.Range("N135").Value = -43302
$ws.Range("H137").Value = 15152618
$ws.Range("I137").Value = 20000764
$ws.Range("J137").Value = 2163.875
$ws.Range("K137").Value = 60002292
$ws.Range("L137").Value = 6491.625
$ws.Range("M137").Value = -59999742
$ws.Range("N137").Value = -11591.625
$ws.Range("H138").Value = 3844.0225
$ws.Range("I138").Value = 1220.44
$ws.Range("J138").Value = 4868.8594
$ws.Range("K138").Value = 3661.32
$ws.Range("L138").Value = 14606.5782
$ws.Range("M138").Value = 1478.68
$ws.Range("N138").Value = -24886.5782

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2004.8667
$ws.Range("J61").Value = 2680.3635
$ws.Range("L61").Value = 2680.3635
$ws.Range("N61").Value = -3104.3635
$ws.Range("H110").Value = 1380.875
$ws.Range("I110").Value = 593
$ws.Range("J110").Value = 2483.9
$ws.Range("K110").Value = 593
$ws.Range("L110").Value = 2483.9
$ws.Range("M110").Value = 1452
$ws.Range("N110").Value = -6573.9
$ws.Range("H122").Value = 5953428
$ws.Range("I122").Value = 6579789
$ws.Range("K122").Value = 19739367
$ws.Range("M122").Value = -19736917
$ws.Range("H136").Value = 2004.8667
$ws.Range("J136").Value = 2680.3635
$ws.Range("L136").Value = 8041.0905
$ws.Range("N136").Value = -13141.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1043.091
$ws.Range("I94").Value = 728.1724
$ws.Range("J94").Value = 1651.9333
$ws.Range("K94").Value = 728.1724
$ws.Range("L94").Value = 1651.9333
$ws.Range("M94").Value = -277.1724
$ws.Range("N94").Value = -2553.9333
$ws.Range("H105").Value = 1895763.5
$ws.Range("I105").Value = 2274316.2
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2274316.2
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -2272569.2
$ws.Range("N105").Value = -6494
$ws.Range("H134").Value = 3290741
$ws.Range("I134").Value = 4238219
$ws.Range("J134").Value = 2435.9412
$ws.Range("K134").Value = 12714657
$ws.Range("L134").Value = 7307.823600000001
$ws.Range("M134").Value = -12712122
$ws.Range("N134").Value = -12377.8236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 759.4286
$ws.Range("I22").Value = 495.25
$ws.Range("J22").Value = 865.1
$ws.Range("K22").Value = 495.25
$ws.Range("L22").Value = 865.1
$ws.Range("M22").Value = -145.25
$ws.Range("N22").Value = -1565.1
$ws.Range("H122").Value = 22602.4
$ws.Range("I122").Value = 34337.332
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 103011.996
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -100561.996
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 657.2857
$ws.Range("I5").Value = 433.53333
$ws.Range("K5").Value = 1300.59999
$ws.Range("M5").Value = -1188.59999
$ws.Range("H129").Value = 1226
$ws.Range("I129").Value = 1043.3334
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 3130.0002
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 1869.9998
$ws.Range("N129").Value = -14500
$ws.Range("H135").Value = 657.2857
$ws.Range("I135").Value = 433.53333
$ws.Range("K135").Value = 3901.79997
$ws.Range("M135").Value = -1366.79997
$ws.Range("H136").Value = 3588.5715
$ws.Range("I136").Value = 3190
$ws.Range("J136").Value = 3748
$ws.Range("K136").Value = 9570
$ws.Range("L136").Value = 11244
$ws.Range("M136").Value = -4470
$ws.Range("N136").Value = -21444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 925.7273
$ws.Range("I46").Value = 672.75
$ws.Range("J46").Value = 1070.2858
$ws.Range("K46").Value = 672.75
$ws.Range("L46").Value = 1070.2858
$ws.Range("M46").Value = -484.75
$ws.Range("N46").Value = -1446.2858
$ws.Range("H55").Value = 158.75
$ws.Range("I55").Value = 57.57143
$ws.Range("J55").Value = 213.23077
$ws.Range("K55").Value = 57.57143
$ws.Range("L55").Value = 213.23077
$ws.Range("M55").Value = 115.42857
$ws.Range("N55").Value = -559.23077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 4865.2
$ws.Range("I45").Value = 1800
$ws.Range("J45").Value = 5631.5
$ws.Range("K45").Value = 1800
$ws.Range("L45").Value = 5631.5
$ws.Range("M45").Value = -1309
$ws.Range("N45").Value = -6613.5
$ws.Range("H74").Value = 8350
$ws.Range("I74").Value = 7800
$ws.Range("J74").Value = 8900
$ws.Range("K74").Value = 7800
$ws.Range("L74").Value = 8900
$ws.Range("M74").Value = -6864
$ws.Range("N74").Value = -10772
$ws.Range("H77").Value = 8350
$ws.Range("I77").Value = 7800
$ws.Range("J77").Value = 8900
$ws.Range("K77").Value = 23400
$ws.Range("L77").Value = 26700
$ws.Range("M77").Value = -18720
$ws.Range("N77").Value = -36060
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H114").Value = 36841.375
$ws.Range("J114").Value = 36841.375
$ws.Range("L114").Value = 36841.375
$ws.Range("N114").Value = -45519.375
$ws.Range("H126").Value = 1484.9375
$ws.Range("I126").Value = 975.9
$ws.Range("J126").Value = 2333.3333
$ws.Range("K126").Value = 2927.7
$ws.Range("L126").Value = 6999.999899999999
$ws.Range("M126").Value = -457.6999999999998
$ws.Range("N126").Value = -11939.9999
$ws.Range("H132").Value = 1397.55
$ws.Range("I132").Value = 1131.3636
$ws.Range("K132").Value = 3394.0908
$ws.Range("M132").Value = -864.0907999999999
